# "Actualización: Base de datos y plantilla Word actualizadas"
# The underlying data update removed one auditee record (RAMIREZ SOSA, JOSE
# GUILLERMO, RFC RASG660608VB3) from the "SUP 29D" database sheet. Deleting
# the entire worksheet row shifts all subsequent records up by one row,
# which matches the rest of the observed changes (dimension, row heights
# carried with their records, shared-string renumbering, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SUP 29D")

# Remove the row for RAMIREZ / SOSA / JOSE GUILLERMO / RASG660608VB3
$ws.Rows("39:39").Delete()

# Refresh the remembered sort range/condition so it covers the new extent
# of the table (previously D2:G44 / D1:D44, now one row shorter).
$sortRange = $ws.Range("D2:G43")
$sortKey = $ws.Range("D1:D43")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Shrink the hidden _FilterDatabase defined name to match the smaller table
$filterName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterName.RefersTo = "='SUP 29D'!`$A`$1:`$G`$38"

# Leave the cursor on the row that now holds the data formerly in row 40
$excel.Goto($ws.Range("A39:XFD39"))
